$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.981.46"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.746.76"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "233.58"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5186"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").Value = "0.2820"
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("D9").Value = "39.52"
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("D10").Value = "0.06130"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "1.755.91"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "0.07012"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "15.40"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "0.6427"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").Value = "4.522"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "77.43"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "25.997.51"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "0.000006609"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "1.978.48"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "4.150"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Value = "8.644"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "5.140"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "139.23"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").Value = "1.505"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.11"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.825"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "102.61"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "0.08271"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "3.665"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "3.434"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "0.04478"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "2.615"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "0.9906"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "0.6157"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "0.01592"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").Value = "1.934"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "100.68"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "0.3859"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "5.068"
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("D46").Value = "0.05466"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  +5.70%  "
$ws.Range("D48").Value = "0.1125"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").Value = "53.02"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "30.06"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "7.621"
$ws.Range("E51").Value = "  +1.68%  "
